# Applies odds value updates to rows 2, 5, 6 and 7 (Sheet1)
# as described by the source diff for Jogos_da_Semana_FlashScore_2024-10-31.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 5.25
$ws.Range("H2").Value = 3.8
$ws.Range("I2").Value = 1.67
$ws.Range("Q2").Value = 2.05
$ws.Range("R2").Value = 1.85
$ws.Range("W2").Value = 13
$ws.Range("X2").Value = 26
$ws.Range("AK2").Value = 12
$ws.Range("AS2").Value = 301
$ws.Range("AX2").Value = 8.5

# Row 5
$ws.Range("G5").Value = 2.4
$ws.Range("H5").Value = 2.6
$ws.Range("I5").Value = 3.6
$ws.Range("J5").Value = 3.15
$ws.Range("K5").Value = 1.78
$ws.Range("M5").Value = 1.16
$ws.Range("N5").Value = 4.55
$ws.Range("O5").Value = 1.65
$ws.Range("P5").Value = 2.15
$ws.Range("Q5").Value = 2.9
$ws.Range("R5").Value = 1.36
$ws.Range("S5").Value = 1.65
$ws.Range("T5").Value = 2.15
$ws.Range("U5").Value = 2.25
$ws.Range("V5").Value = 1.57
$ws.Range("W5").Value = 5.4
$ws.Range("X5").Value = 10
$ws.Range("Y5").Value = 10
$ws.Range("Z5").Value = 27
$ws.Range("AA5").Value = 27
$ws.Range("AB5").Value = 50
$ws.Range("AC5").Value = 4.55
$ws.Range("AE5").Value = 19.5
$ws.Range("AF5").Value = 150
$ws.Range("AH5").Value = 7
$ws.Range("AI5").Value = 17
$ws.Range("AJ5").Value = 13.5
$ws.Range("AK5").Value = 60
$ws.Range("AL5").Value = 50
$ws.Range("AM5").Value = 70
$ws.Range("AN5").Value = 4
$ws.Range("AO5").Value = 14
$ws.Range("AP5").Value = 27
$ws.Range("AQ5").Value = 70
$ws.Range("AR5").Value = 150
$ws.Range("AS5").Value = 500
$ws.Range("AT5").Value = 2.1
$ws.Range("AU5").Value = 7.9
$ws.Range("AV5").Value = 100
$ws.Range("AW5").Value = 5.2
$ws.Range("AX5").Value = 24
$ws.Range("AY5").Value = 37
$ws.Range("BA5").Value = 250

# Row 6
$ws.Range("G6").Value = 2.12
$ws.Range("H6").Value = 4.1
$ws.Range("I6").Value = 2.75
$ws.Range("J6").Value = 2.52
$ws.Range("K6").Value = 2.6
$ws.Range("L6").Value = 3.05
$ws.Range("O6").Value = 1.1
$ws.Range("P6").Value = 6
$ws.Range("Q6").Value = 1.32
$ws.Range("R6").Value = 3.1
$ws.Range("S6").Value = 1.19
$ws.Range("T6").Value = 4.15
$ws.Range("U6").Value = 1.32
$ws.Range("V6").Value = 3.1
$ws.Range("W6").Value = 17
$ws.Range("X6").Value = 17
$ws.Range("Y6").Value = 9.75
$ws.Range("Z6").Value = 24
$ws.Range("AB6").Value = 16
$ws.Range("AC6").Value = 10.5
$ws.Range("AD6").Value = 9.75
$ws.Range("AF6").Value = 25
$ws.Range("AG6").Value = 100
$ws.Range("AH6").Value = 20
$ws.Range("AI6").Value = 23
$ws.Range("AJ6").Value = 11.5
$ws.Range("AK6").Value = 37
$ws.Range("AM6").Value = 18
$ws.Range("AN6").Value = 5
$ws.Range("AP6").Value = 12.5
$ws.Range("AQ6").Value = 32
$ws.Range("AR6").Value = 37
$ws.Range("AS6").Value = 90
$ws.Range("AT6").Value = 4.15
$ws.Range("AV6").Value = 29
$ws.Range("AW6").Value = 5.7
$ws.Range("AY6").Value = 14
$ws.Range("BA6").Value = 50
$ws.Range("BB6").Value = 100
$ws.Range("BC6").Value = 300

# Row 7
$ws.Range("G7").Value = 2.35
$ws.Range("I7").Value = 2.63
$ws.Range("J7").Value = 2.75
$ws.Range("S7").Value = 1.25
$ws.Range("T7").Value = 3.75
$ws.Range("X7").Value = 17
$ws.Range("AA7").Value = 17
$ws.Range("AH7").Value = 15
$ws.Range("AJ7").Value = 11
$ws.Range("AO7").Value = 12
$ws.Range("AT7").Value = 3.75

Write-Output "Applied odds updates to rows 2, 5, 6, 7"
